$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SocialFeedID" column entirely (was column A), shifting
# Date/Link/CompanyID/CompanyName one column to the left.
$null = $ws.Range("A1:A1048576").Delete()

# Add a new data row (row 2): a date, a hyperlinked URL, and a company id.
$ws.Range("A2").Value = 35864
$null = $ws.Hyperlinks.Add($ws.Range("B2"), "https://pnq.co.in/")
$ws.Range("C2").Value = "DEMC"

# Update the active selection to match the saved workbook state.
$null = $ws.Range("B13").Select()
